# Pushing to products and product-pool is now working, on delete will
# create an alert now.
#
#  - TODO sheet: "IP-Adressen loggen und deren Tätigkeit, am besten auch
#    User-Ids" (row 27) moves from "offen" to "done"; add a new TODO item
#    "console.logs entfernen" (row 30, status "offen").
#  - TODO CMS sheet: "Eingabe testen" (row 25) moves from "offen" to
#    "done"; add a new TODO item "console.logs entfernen" (row 26, status
#    "offen").

$wb = $excel.ActiveWorkbook

# --- Sheet "TODO" ---
$wsTodo = $wb.Worksheets.Item("TODO")

# Row 27: "offen" -> "done" (copy the existing "done" cell format so the
# shared green style is reused instead of minting a new one).
$wsTodo.Range("B2").Copy()
$wsTodo.Range("B27").PasteSpecial(-4122)  # xlPasteFormats
$wsTodo.Range("B27").Value = "done"

# New row 30: "console.logs entfernen" / "offen".
$wsTodo.Range("A30").Value = "console.logs entfernen"
$wsTodo.Range("B4").Copy()
$wsTodo.Range("B30").PasteSpecial(-4122)  # xlPasteFormats ("offen" style)
$wsTodo.Range("B30").Value = "offen"

$wsTodo.Range("A35").Select()

# --- Sheet "TODO CMS" ---
$wsCms = $wb.Worksheets.Item("TODO CMS")

# Row 25: "offen" -> "done".
$wsCms.Range("B2").Copy()
$wsCms.Range("B25").PasteSpecial(-4122)  # xlPasteFormats
$wsCms.Range("B25").Value = "done"

# New row 26: "console.logs entfernen" / "offen".
$wsCms.Range("A26").Value = "console.logs entfernen"
$wsCms.Range("B24").Copy()
$wsCms.Range("B26").PasteSpecial(-4122)  # xlPasteFormats ("offen" style)
$wsCms.Range("B26").Value = "offen"

$wsCms.Activate()
$wsCms.Range("B15").Select()

$excel.CutCopyMode = $false
